# Add a new signup row (row 6) to the "User Signups" sheet, mirroring the
# layout used by the existing short-form rows (4 and 5): Timestamp, Full
# Name, Email, two extra fields, Registration IP, Status (7 columns, A:G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "2025-10-27T11:02:02.845Z"
$ws.Range("B6").Value = "Patrick Sharma"
$ws.Range("C6").Value = "patricksharma1234@gmail.com"
$ws.Range("D6").Value = "Ca23m6Na"
$ws.Range("E6").Value = "HBhjbs234"
$ws.Range("F6").Value = "::1"
$ws.Range("G6").Value = "Active"

# Match the formatting of the previous short-form row (row 4) so the new
# row picks up the same fill/border style used for alternating rows.
$srcFormat = $ws.Range("A4:G4")
$dstFormat = $ws.Range("A6:G6")
$srcFormat.Copy()
$dstFormat.PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = 0
